$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row for Pepe Lopez (row 7). Columns D and G look numeric
# ("12121212", "44300") but must stay text, matching the rest of the
# sheet (e.g. NumeroTelefonico/CodigoPostal columns in earlier rows),
# so force those two cells to Text format before writing them.
$ws.Range("D7").NumberFormat = "@"
$ws.Range("G7").NumberFormat = "@"

$ws.Range("A7").Value = "Pepe"
$ws.Range("B7").Value = "Lopez"
$ws.Range("C7").Value = "PepeL"
$ws.Range("D7").Value = "12121212"
$ws.Range("E7").Value = "pepel@ejemplo.com"
$ws.Range("F7").Value = "fkjdfhaskjfhashfhasfadfhjasdfhads"
$ws.Range("G7").Value = "44300"
$ws.Range("H7").Value = "Sin Adquirir"
